# chore: adapt column header formatting to respective input file names
#
#   *_old  -> *_FV2410   (the "before" / FV2410 format version columns)
#   *_new  -> *_FV2504   (the "after"  / FV2504 format version columns)
#
# Additionally expose the data range as a proper Excel Table (so the
# header row gets filters + structured references) and freeze the
# header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row -------------------------------------------
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
  $headerCell = $ws.Cells.Item(1, $c)
  $headerText = $headerCell.Value()
  if ($headerText -like "*_old") {
    $headerCell.Value = $headerText.Substring(0, $headerText.Length - 4) + "_FV2410"
  } elseif ($headerText -like "*_new") {
    $headerCell.Value = $headerText.Substring(0, $headerText.Length - 4) + "_FV2504"
  }
}

# --- 2. Turn the data range into a native Excel Table -------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row --------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
